# Change to 75 wars; score formula update (less points for total victories)
# Update member stats in the "Info" sheet to reflect the new scoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Row 4 (Pipa)
$ws.Range("O4").Value = 662.0
$ws.Range("P4").Value = 568.0
$ws.Range("Q4").Value = 99727.0

# Row 17 (bascenso)
$ws.Range("P17").Value = 757.0

# Row 18 (thunder)
$ws.Range("P18").Value = 717.0

# Row 24 (Ribiti)
$ws.Range("D24").Value = 5342.0
$ws.Range("G24").Value = 13524.0
$ws.Range("H24").Value = 43061.0

# Row 41 (Savler)
$ws.Range("F41").Value = 4471.0
$ws.Range("H41").Value = 13473.0
$ws.Range("I41").Value = 3969.0

# Row 49 (Rafa)
$ws.Range("O49").Value = 766.0
$ws.Range("Q49").Value = 37604.0
